$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 28
$ws.Range("F2").Value = 452.53
$ws.Range("G2").Value = 10.69
$ws.Range("H2").Value = 18.94
$ws.Range("I2").Value = 1.27
$ws.Range("J2").Value = 0.11
$ws.Range("K2").Value = 4.75
$ws.Range("L2").Value = 0.72
$ws.Range("M2").Value = 0.06
$ws.Range("N2").Value = 14.9
$ws.Range("O2").Value = 1.09
$ws.Range("P2").Value = 0.09
$ws.Range("Q2").Value = 66.54000000000001
$ws.Range("R2").Value = 12.37
$ws.Range("S2").Value = 1.05
$ws.Range("T2").Value = 4.74
$ws.Range("U2").Value = 0.78
$ws.Range("V2").Value = 0.07000000000000001
$ws.Range("W2").Value = 324.02
$ws.Range("X2").Value = 27.51
$ws.Range("Y2").Value = 2.34
$ws.Range("Z2").Value = 23.53
$ws.Range("AA2").Value = 1.66
$ws.Range("AB2").Value = 0.14
$ws.Range("AC2").Value = 19.72
$ws.Range("AD2").Value = 1.51
$ws.Range("AE2").Value = 0.13
$ws.Range("AF2").Value = 20.78
$ws.Range("AG2").Value = 1.59
$ws.Range("AH2").Value = 0.13
$ws.Range("AI2").Value = 24.73
$ws.Range("AJ2").Value = 0.58
$ws.Range("AK2").Value = 0.05
